$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Remove existing hyperlinks (avoids stale/misaligned refs after rewrite) ---
$guard = 0
while ($ws.Hyperlinks.Count -gt 0 -and $guard -lt 50) {
    foreach ($hl in $ws.Hyperlinks) {
        $hl.Delete()
    }
    $guard = $guard + 1
}

# --- Column H width: 14 -> 16 chars (ColumnWidth offset ~0.83 vs stored char width) ---
$ws.Range("H1").ColumnWidth = 15.17

# --- Rewrite data rows 2-13 with final content ---
$ws.Range("A2").Value = "2026-02-10 13:07:34"
$ws.Range("B2").Value = "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5423720"
$ws.Range("G2").Value = 385
$ws.Range("H2").Value = "🔥AI,Ai ◆効率化"

$ws.Range("A3").Value = "2026-02-10 13:07:34"
$ws.Range("B3").Value = "建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5434128"
$ws.Range("G3").Value = 368
$ws.Range("H3").Value = "🔥AI,Ai ◆開発"

$ws.Range("A4").Value = "2026-02-10 13:07:34"
$ws.Range("B4").Value = "企業のMicrosoft Copilot導入・活用支援AIコンサルタント募集(研修講師・メンター)"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5434363"
$ws.Range("G4").Value = 348
$ws.Range("H4").Value = "🔥AI,Ai ◆コンサル"

$ws.Range("A5").Value = "2026-02-10 13:07:34"
$ws.Range("B5").Value = "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5427956"
$ws.Range("G5").Value = 310
$ws.Range("H5").Value = "🔥AI,Ai"

$ws.Range("A6").Value = "2026-02-10 13:07:34"
$ws.Range("B6").Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Range("G6").Value = 243
$ws.Range("H6").Value = "🔥API ◆ツール"

$ws.Range("A7").Value = "2026-02-10 13:07:34"
$ws.Range("B7").Value = "【急募】新聞記事PDFをCSV・Excel化するPythonプログラム作成依頼"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5489128"
$ws.Range("G7").Value = 198
$ws.Range("H7").Value = "🔥Python"

$ws.Range("A8").Value = "2026-02-10 13:07:34"
$ws.Range("B8").Value = "初回 自動車販売・整備業の管理システム開発"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5489393"
$ws.Range("G8").Value = 153
$ws.Range("H8").Value = "◆開発,システム開発 ◇管理"

$ws.Range("A9").Value = "2026-02-10 13:07:34"
$ws.Range("B9").Value = "【Unity/XRエンジニア募集】製造業DX支援!既存システムと連携するXRアプリ開発"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5454210"
$ws.Range("G9").Value = 108
$ws.Range("H9").Value = "◆開発 ◇アプリ"

$ws.Range("A10").Value = "2026-02-10 13:07:34"
$ws.Range("B10").Value = "【外国人大歓迎】【急募】ECツールの保守・バグ修正・機能追加エンジニア募集"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5489500"
$ws.Range("G10").Value = 68
$ws.Range("H10").Value = "◆ツール"

$ws.Range("A11").Value = "2026-02-10 13:07:34"
$ws.Range("B11").Value = "スプレッドシート(Apps Script)で作業時間をボタン1つで計測・集計できる仕組みの開発"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5488743"
$ws.Range("G11").Value = 68
$ws.Range("H11").Value = "◆開発"

$ws.Range("A12").Value = "2026-02-10 13:07:34"
$ws.Range("B12").Value = "【農機具管理】顧客指定で保有機情報を見れるシステム構築依頼"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5489112"
$ws.Range("G12").Value = 45
$ws.Range("H12").Value = "◇管理"

$ws.Range("A13").Value = "2026-02-10 13:07:34"
$ws.Range("B13").Value = "【急募】ECサイト(WooCommerce)の決済・配送ロジックテスト、デバッグ検証依頼"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5489409"
$ws.Range("G13").Value = 33
$ws.Range("H13").Value = "◇サイト"

# --- Re-add hyperlinks for F2:F13, then force Hyperlink style reuse (keeps style index 1) ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5423720")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5434128")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5434363")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5427956")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5489128")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5489393")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5454210")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5489500")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5488743")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5489112")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5489409")
$ws.Range("F2:F13").Style = "Hyperlink"

Write-Output "done"
